$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the ASV_ID / Species_name / Common_name / Category values between
# row 42 (Homo sapiens / Human) and row 43 (unassigned), leaving the
# ASV_sum (E) and ASV_rank (F) columns untouched.

$a42 = $ws.Range("A42").Value2
$b42 = $ws.Range("B42").Value2
$c42 = $ws.Range("C42").Value2
$d42 = $ws.Range("D42").Value2

$a43 = $ws.Range("A43").Value2
$b43 = $ws.Range("B43").Value2
$c43 = $ws.Range("C43").Value2
$d43 = $ws.Range("D43").Value2

$ws.Range("A42").Value2 = $a43
$ws.Range("B42").Value2 = $b43
$ws.Range("C42").Value2 = $c43
$ws.Range("D42").Value2 = $d43

$ws.Range("A43").Value2 = $a42
$ws.Range("B43").Value2 = $b42
$ws.Range("C43").Value2 = $c42
$ws.Range("D43").Value2 = $d42
